$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the existing user's CPF (make it a real number, matching the
# rest of the sheet) and fix up the user's name.
$ws.Range("B2").Value = 11122233304
$ws.Range("C2").Value = "Usuario testado"

# Add the new user created on row 3.
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 11122233305
$ws.Range("C3").Value = "Teste"

# Resize column B (cpf) to fit its new, wider numeric content.
$ws.Columns.Item(2).AutoFit() | Out-Null

# Leave the selection on the next empty row, like a user does after
# finishing data entry.
$ws.Rows.Item(4).Select() | Out-Null
